# Adapt column header formatting to respective input file names (#7)
#
# 1. Rename the column headers in row 1:
#      "<name>_old" -> "<name>_FV2210"
#      "<name>_new" -> "<name>_FV2304"
#    (column "diff" in between stays as-is)
# 2. Freeze the header row (freeze panes at A2).
# 3. Turn the used data range A1:U67 into a proper Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerBase = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10): old/"FV2210" headers
$col = 1
foreach ($name in $headerBase) {
    $ws.Cells.Item(1, $col).Value = "$($name)_FV2210"
    $col = $col + 1
}

# Column K (11): "diff" header, unchanged
$col = $col + 1

# Columns L-U (12-21): new/"FV2304" headers
foreach ($name in $headerBase) {
    $ws.Cells.Item(1, $col).Value = "$($name)_FV2304"
    $col = $col + 1
}

# Freeze the header row (row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Convert the data range into an Excel Table / ListObject.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U67"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
